$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("latest")

$ws.Range("B2").Value = -0.3731803712827915
$ws.Range("C2").Value = 2.007670690660352
$ws.Range("D2").Value = 12.87239756141505
$ws.Range("E2").Value = 3.58781236429876
$ws.Range("F2").Value = 3.652324422179951
$ws.Range("G2").Value = 22
$ws.Range("B3").Value = -0.4804827079566977
$ws.Range("C3").Value = 1.161321962968468
$ws.Range("D3").Value = 5.130197931044437
$ws.Range("E3").Value = 2.264994024505238
$ws.Range("F3").Value = 2.268105159228511
$ws.Range("G3").Value = 21
$ws.Range("B4").Value = -0.2697575283137488
$ws.Range("C4").Value = 0.8727440841527818
$ws.Range("D4").Value = 2.711211561858076
$ws.Range("E4").Value = 1.64657570790355
$ws.Range("F4").Value = 1.666525676141245
$ws.Range("G4").Value = 20
$ws.Range("B5").Value = -0.02251358845263601
$ws.Range("C5").Value = 0.8231985852708921
$ws.Range("D5").Value = 1.758851113115051
$ws.Range("E5").Value = 1.326216842418709
$ws.Range("F5").Value = 1.362361935462548
$ws.Range("G5").Value = 19
$ws.Range("B6").Value = -0.1079376821351743
$ws.Range("C6").Value = 0.7903692723160057
$ws.Range("D6").Value = 1.18367369466905
$ws.Range("E6").Value = 1.087967690085073
$ws.Range("F6").Value = 1.113986395682008
$ws.Range("G6").Value = 18
$ws.Range("B7").Value = 0.05847297744281363
$ws.Range("C7").Value = 0.6864530745518149
$ws.Range("D7").Value = 0.8768293907560101
$ws.Range("E7").Value = 0.9363916866119701
$ws.Range("F7").Value = 0.9633267594741898
$ws.Range("G7").Value = 17
$ws.Range("B8").Value = 0.1260832718587093
$ws.Range("C8").Value = 0.5907394480024477
$ws.Range("D8").Value = 0.5291801755352681
$ws.Range("E8").Value = 0.7274477132655432
$ws.Range("F8").Value = 0.7399338234140823
$ws.Range("G8").Value = 16
$ws.Range("B9").Value = 0.189738633461222
$ws.Range("C9").Value = 0.6135846084182608
$ws.Range("D9").Value = 0.6359933741545197
$ws.Range("E9").Value = 0.7974919273287472
$ws.Range("F9").Value = 0.8017787981776971
$ws.Range("G9").Value = 15
$ws.Range("B10").Value = 0.2282491029459698
$ws.Range("C10").Value = 0.5916736138517826
$ws.Range("D10").Value = 0.4825523572072717
$ws.Range("E10").Value = 0.6946598859926142
$ws.Range("F10").Value = 0.6808572570925593
$ws.Range("G10").Value = 14
$ws.Range("B11").Value = 0.2349671465221997
$ws.Range("C11").Value = 0.5153392861830149
$ws.Range("D11").Value = 0.3670491565632086
$ws.Range("E11").Value = 0.6058458191348758
$ws.Range("F11").Value = 0.5812282079097323
$ws.Range("G11").Value = 13
